# Added waitbar to Generate button
#
# The "Generate" export routine now reports a Scaffold *Width* instead of a
# Scaffold *Diameter*, and the exported info sheet is refreshed with the
# generator's current (default) parameter values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (approximate refresh of the auto-sized columns)
$ws.Columns.Item(1).ColumnWidth = 21.60807291666667
$ws.Columns.Item(2).ColumnWidth = 15.514322916666666
$ws.Columns.Item(3).ColumnWidth = 14.334635416666666
$ws.Columns.Item(4).ColumnWidth = 13.787760416666666

# Header row
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Original Geometry"
$ws.Range("C1").Value = "Scaled Geometry"
$ws.Range("D1").Value = "Fitted Geometry"

# Row labels (column A)
$ws.Range("A2").Value = "Scaffold Length (mm)"
$ws.Range("A3").Value = "Scaffold Width (mm)"
$ws.Range("A4").Value = "Pore Length: p (µm)"
$ws.Range("A5").Value = "Pore Width: q (µm)"
$ws.Range("A6").Value = "Pore Area: K (mm²)"
$ws.Range("A7").Value = "Pore Perimeter: P (mm)"
$ws.Range("A8").Value = "Major Pore Angle: A (deg)"
$ws.Range("A9").Value = "Minor Pore Angle: B (deg)"
$ws.Range("A10").Value = "Repeating Cells (X)"
$ws.Range("A11").Value = "Repeating Cells (Y)"

# Data values (Original / Scaled / Fitted geometry columns)
$ws.Range("B2:D2").Value = 1000
$ws.Range("B3:D3").Value = 1000
$ws.Range("B4:D4").Value = 1000
$ws.Range("B5:D5").Value = 1000
$ws.Range("B6:D6").Value = 1
$ws.Range("B7:D7").Value = 4

$ws.Range("B8:D8").ClearContents()
$ws.Range("B9:D9").ClearContents()

$ws.Range("B10:D10").Value = 500
$ws.Range("B11:D11").Value = 500
